$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'57.847.26"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.05%  "

$ws.Range("D3").Value = "'2.288.45"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.09%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'536.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.33%  "

$ws.Range("D6").Value = "'131.22"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.26%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("E8").Value = "  -3.34%  "

$ws.Range("D9").Value = "'2.287.12"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.98%  "

$ws.Range("E10").Value = "  -5.65%  "

$ws.Range("D11").Value = "'5.46"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.84%  "

$ws.Range("E12").Value = "  -0.60%  "

$ws.Range("E13").Value = "  -4.83%  "

$ws.Range("D14").Value = "'23.54"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.70%  "

$ws.Range("D15").Value = "'2.700.44"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.85%  "

$ws.Range("D16").Value = "'57.875.41"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.86%  "

$ws.Range("E17").Value = "  -4.90%  "

$ws.Range("D18").Value = "'2.330.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.39%  "

$ws.Range("D19").Value = "'10.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.82%  "

$ws.Range("E20").Value = "  -6.67%  "

$ws.Range("D21").Value = "'313.84"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.43%  "

$ws.Range("D22").Value = "'6.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.42%  "

$ws.Range("D23").Value = "'0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.01%  "

$ws.Range("D24").Value = "'62.99"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.40%  "

$ws.Range("E25").Value = "  -4.89%  "

$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.01%  "

$ws.Range("D27").Value = "'7.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.79%  "

$ws.Range("E28").Value = "  -5.61%  "

$ws.Range("D29").Value = "'1.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.41%  "

$ws.Range("D30").Value = "'169.68"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.53%  "

$ws.Range("D31").Value = "'0.0₃0720"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.89%  "

$ws.Range("E32").Value = "  -0.49%  "

$ws.Range("D33").Value = "'5.73"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.58%  "

$ws.Range("E34").Value = "  -5.82%  "

$ws.Range("E35").Value = "  -0.02%  "

$ws.Range("D36").Value = "'17.67"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.86%  "

$ws.Range("E37").Value = "  -0.01%  "

$ws.Range("E38").Value = "  -7.81%  "

$ws.Range("D39").Value = "'3.90"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.74%  "

$ws.Range("E40").Value = "  -1.88%  "

$ws.Range("E41").Value = "  -7.29%  "

$ws.Range("D42").Value = "'141.63"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.17%  "

$ws.Range("D43").Value = "'288.30"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -10.94%  "

$ws.Range("E44").Value = "  -4.59%  "

$ws.Range("D46").Value = "'0.0498"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.72%  "

$ws.Range("D47").Value = "'0.553"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.72%  "

$ws.Range("D48").Value = "'18.21"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -8.34%  "

$ws.Range("D49").Value = "'0.0210"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.21%  "

$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D50").Value = "'10.96"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.88%  "

$ws.Range("B51").Value = "ZEEBU"
$ws.Range("C51").Value = "https://coinranking.com/coin/B5-YKN_zB+zeebu-zbu"
$ws.Range("D51").Value = "'4.65"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.35%  "
